$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark ---
# Originally it sits right before "Delfina da Silva Manente" and its (accidental)
# closing tag sits in the middle of the "ministraram cursos diversos" sentence.
# The edit removes it from both of those spots and re-creates it, collapsed,
# right after the "B" of "Brasileira" (second paragraph), once the leading
# space before "Brasileira" has been dropped.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Remove the stray leading space before "Brasileira, 38 anos" ---
$d.Content.Find.Execute(" Brasileira, 38 anos ", $true, $false, $false, $false, $false, $true, 1, $false, "Brasileira, 38 anos ", 2) | Out-Null

# Re-insert the "_GoBack" bookmark, collapsed, right after the initial "B".
$rng = $d.Content
$rng.Find.Execute("Brasileira") | Out-Null
$bmRange = $d.Range($rng.Start + 1, $rng.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 3. Merge runs that were needlessly split (no textual change, just
#        collapsing sibling runs that share identical formatting) ---

$d.Content.Find.Execute("Buffet Vila dos brutos 03/07/2019 à 17/03/2020 ", $true, $false, $false, $false, $false, $true, 1, $false, "Buffet Vila dos brutos 03/07/2019 à 17/03/2020 ", 2) | Out-Null

$d.Content.Find.Execute("Preparação de pratos quentes e frios, café da manhã, ", $true, $false, $false, $false, $false, $true, 1, $false, "Preparação de pratos quentes e frios, café da manhã, ", 2) | Out-Null

$d.Content.Find.Execute(" lance em eventos (cozinha) Clube Royal Five e outros, de 01/2008 a 04/2009 ", $true, $false, $false, $false, $false, $true, 1, $false, " lance em eventos (cozinha) Clube Royal Five e outros, de 01/2008 a 04/2009 ", 2) | Out-Null

$d.Content.Find.Execute("Tarefas relacionadas à cozinha: ", $true, $false, $false, $false, $false, $true, 1, $false, "Tarefas relacionadas à cozinha: ", 2) | Out-Null

$d.Content.Find.Execute("Preparação de jantares especiais, café da manhã e ", $true, $false, $false, $false, $false, $true, 1, $false, "Preparação de jantares especiais, café da manhã e ", 2) | Out-Null

$d.Content.Find.Execute("Montagem de pratos quentes e saladas variadas; ", $true, $false, $false, $false, $false, $true, 1, $false, "Montagem de pratos quentes e saladas variadas; ", 2) | Out-Null

$d.Content.Find.Execute("Preparação de alimentos e coquetéis; ", $true, $false, $false, $false, $false, $true, 1, $false, "Preparação de alimentos e coquetéis; ", 2) | Out-Null

$d.Content.Find.Execute("Tarefas administrativas: ", $true, $false, $false, $false, $false, $true, 1, $false, "Tarefas administrativas: ", 2) | Out-Null

$d.Content.Find.Execute("Pagamentos (internos e em bancos); ", $true, $false, $false, $false, $false, $true, 1, $false, "Pagamentos (internos e em bancos); ", 2) | Out-Null

$d.Content.Find.Execute(", Carla Pernambuco, Luiz Cintra e outros) ministraram cursos diversos, como ", $true, $false, $false, $false, $false, $true, 1, $false, ", Carla Pernambuco, Luiz Cintra e outros) ministraram cursos diversos, como ", 2) | Out-Null
